$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1091.641
$ws.Range("I15").Value = 1091.641
$ws.Range("K15").Value = 3274.923
$ws.Range("M15").Value = -3105.923

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 758552
$ws.Range("I80").Value = 947897.4
$ws.Range("J80").Value = 1170.5
$ws.Range("K80").Value = 2843692.2
$ws.Range("L80").Value = 3511.5
$ws.Range("M80").Value = -2842694.2
$ws.Range("N80").Value = -5507.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 758552
$ws.Range("I83").Value = 947897.4
$ws.Range("J83").Value = 1170.5
$ws.Range("K83").Value = 8531076.6
$ws.Range("L83").Value = 10534.5
$ws.Range("M83").Value = -8526084.6
$ws.Range("N83").Value = -20518.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 24361268
$ws.Range("I116").Value = 25112032
$ws.Range("J116").Value = 23825006
$ws.Range("K116").Value = 25112032
$ws.Range("L116").Value = 23825006
$ws.Range("M116").Value = -25108590
$ws.Range("N116").Value = -23831890

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3800.0312
$ws.Range("I32").Value = 4137.3335
$ws.Range("J32").Value = 1978.6
$ws.Range("K32").Value = 4137.3335
$ws.Range("L32").Value = 1978.6
$ws.Range("M32").Value = -3850.3335
$ws.Range("N32").Value = -2552.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2418.0435
$ws.Range("I74").Value = 2255.2727
$ws.Range("K74").Value = 2255.2727
$ws.Range("M74").Value = -1381.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2418.0435
$ws.Range("I77").Value = 2255.2727
$ws.Range("K77").Value = 11276.3635
$ws.Range("M77").Value = -6908.363499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19443.318
$ws.Range("I20").Value = 35979
$ws.Range("J20").Value = 2907.6365
$ws.Range("K20").Value = 35979
$ws.Range("L20").Value = 2907.6365
$ws.Range("M20").Value = -35732
$ws.Range("N20").Value = -3401.6365

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 19312340
$ws.Range("I107").Value = 84833.44
$ws.Range("K107").Value = 84833.44
$ws.Range("M107").Value = -82913.44

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3749.1667
$ws.Range("J31").Value = 4365.3125
$ws.Range("L31").Value = 4365.3125
$ws.Range("N31").Value = -4955.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3749.1667
$ws.Range("J34").Value = 4365.3125
$ws.Range("L34").Value = 4365.3125
$ws.Range("N34").Value = -4769.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2421.1538
$ws.Range("I58").Value = 1436.7
$ws.Range("K58").Value = 1436.7
$ws.Range("M58").Value = -1233.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 90917640
$ws.Range("I62").Value = 7833.3335
$ws.Range("K62").Value = 7833.3335
$ws.Range("M62").Value = -7209.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 90917640
$ws.Range("I65").Value = 7833.3335
$ws.Range("K65").Value = 39166.6675
$ws.Range("M65").Value = -36046.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1249.6666
$ws.Range("I94").Value = 1059.2
$ws.Range("J94").Value = 1385.7142
$ws.Range("K94").Value = 1059.2
$ws.Range("L94").Value = 1385.7142
$ws.Range("M94").Value = -608.2
$ws.Range("N94").Value = -2287.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2304
$ws.Range("I134").Value = 1861.5714
$ws.Range("K134").Value = 5584.7142
$ws.Range("M134").Value = -3049.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2421.1538
$ws.Range("I136").Value = 1436.7
$ws.Range("K136").Value = 4310.1
$ws.Range("M136").Value = -1760.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 611.4737
$ws.Range("I2").Value = 881.1667
$ws.Range("J2").Value = 149.14285
$ws.Range("K2").Value = 5287.0002
$ws.Range("L2").Value = 894.8571000000001
$ws.Range("M2").Value = -5174.0002
$ws.Range("N2").Value = -1120.8571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9537273
$ws.Range("I4").Value = 14105258
$ws.Range("J4").Value = 4087.1304
$ws.Range("K4").Value = 42315774
$ws.Range("L4").Value = 12261.3912
$ws.Range("M4").Value = -42315662
$ws.Range("N4").Value = -12485.3912

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I17").Value = 416.5
$ws.Range("J17").Value = 867.1429000000001
$ws.Range("K17").Value = 1249.5
$ws.Range("L17").Value = 2601.4287
$ws.Range("M17").Value = -1080.5
$ws.Range("N17").Value = -2939.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2772
$ws.Range("J34").Value = 2951.2727
$ws.Range("L34").Value = 8853.8181
$ws.Range("N34").Value = -9021.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4163.091
$ws.Range("J39").Value = 4479.5
$ws.Range("L39").Value = 13438.5
$ws.Range("N39").Value = -14026.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 896.1667
$ws.Range("I55").Value = 515.4
$ws.Range("J55").Value = 2800
$ws.Range("K55").Value = 1546.2
$ws.Range("L55").Value = 8400
$ws.Range("M55").Value = -1369.2
$ws.Range("N55").Value = -8754

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6847
$ws.Range("I56").Value = 6847
$ws.Range("K56").Value = 6847
$ws.Range("M56").Value = -6317

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 375
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2431.6
$ws.Range("I132").Value = 1499
$ws.Range("K132").Value = 13491
$ws.Range("M132").Value = -10961

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5563.8
$ws.Range("J80").Value = 3956
$ws.Range("L80").Value = 3956
$ws.Range("N80").Value = -5952

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5563.8
$ws.Range("J83").Value = 3956
$ws.Range("L83").Value = 19780
$ws.Range("N83").Value = -29764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3577.7693
$ws.Range("I113").Value = 2319.6
$ws.Range("K113").Value = 2319.6
$ws.Range("M113").Value = -149.5999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 41668772
$ws.Range("I122").Value = 2324.2
$ws.Range("J122").Value = 250001000
$ws.Range("K122").Value = 6972.599999999999
$ws.Range("L122").Value = 750003000
$ws.Range("M122").Value = -4522.599999999999
$ws.Range("N122").Value = -750007900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4320.1787
$ws.Range("I132").Value = 4052.6
$ws.Range("K132").Value = 12157.8
$ws.Range("M132").Value = -9627.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 21996.375
$ws.Range("I61").Value = 22251
$ws.Range("J61").Value = 21741.75
$ws.Range("K61").Value = 22251
$ws.Range("L61").Value = 21741.75
$ws.Range("M61").Value = -22049
$ws.Range("N61").Value = -22145.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 21996.375
$ws.Range("I113").Value = 22251
$ws.Range("J113").Value = 21741.75
$ws.Range("K113").Value = 22251
$ws.Range("L113").Value = 21741.75
$ws.Range("M113").Value = -20081
$ws.Range("N113").Value = -26081.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4369.4287
$ws.Range("I122").Value = 3574.6667
$ws.Range("J122").Value = 5800
$ws.Range("K122").Value = 10724.0001
$ws.Range("L122").Value = 17400
$ws.Range("M122").Value = -8274.000100000001
$ws.Range("N122").Value = -22300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 94875
$ws.Range("J133").Value = 94875
$ws.Range("L133").Value = 94875
$ws.Range("N133").Value = -99935

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 640.55554
$ws.Range("I107").Value = 640.55554
$ws.Range("K107").Value = 1921.66662
$ws.Range("M107").Value = -1.666619999999966

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5877.75
$ws.Range("I113").Value = 1910.8
$ws.Range("K113").Value = 5732.4
$ws.Range("M113").Value = -3562.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3355.1052
$ws.Range("I132").Value = 3009.5312
$ws.Range("K132").Value = 9028.5936
$ws.Range("M132").Value = -6498.5936

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 59998
$ws.Range("J135").Value = 59998
$ws.Range("L135").Value = 59998
$ws.Range("N135").Value = -70138
